$wb = $excel.ActiveWorkbook

$sheetNames = @("Kayitlar", "Erdemli")

$values = @(
    "2987",
    "2025-09-10",
    "Erdemli",
    "1",
    "PAYDAŞ KURUM TALEP",
    "CEMAL TİMUROĞLU (K.Teknisyeni), ALİ BAŞKURT (Kontrol Memuru), SERDAR ARSLAN (Tekniker), EMİNE ALANLI KIRCILI (K.Mühendisi)"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $row = 53

    for ($col = 1; $col -le 6; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.Value = "'" + $values[$col - 1]
        $cell.Style = "Normal"
    }
}
